$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "antenne"
$ws.Range("K2").Value = "MONTREUIL"
$ws.Range("K3").Select()
